$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0

# Row 3
$ws.Range("U3").Value = 0.23
$ws.Range("V3").Value = 0.45
$ws.Range("X3").Value = 0.08
$ws.Range("Y3").Value = 0.04
$ws.Range("Z3").Value = 0.01
$ws.Range("AA3").Value = 0
$ws.Range("AC3").Value = 0

# Row 4
$ws.Range("U4").Value = 0.93
$ws.Range("V4").Value = 0.07000000000000001

# Row 5
$ws.Range("U5").Value = 0.11
$ws.Range("V5").Value = 0.18
$ws.Range("W5").Value = 0.2
$ws.Range("X5").Value = 0.21
$ws.Range("Y5").Value = 0.09
$ws.Range("AB5").Value = 0.03
